$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.523.61"
$ws.Range("D3").Value = "2.108.74"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.71"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5240"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4545"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.39"
$ws.Range("E9").Value = "  +15.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08999"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.52"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "2.107.26"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.782"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.852"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06624"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.34"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.311"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "30.577.70"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.37"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("E25").Value = "  +3.80%  "
$ws.Range("D26").Value = "2.354.96"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.40"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.577"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.36"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.95"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.202"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.663"
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.164"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.56"
$ws.Range("E36").Value = "  +11.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02582"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06809"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.78"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2290"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6928"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.259"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.389"
$ws.Range("E44").Value = "  +6.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6418"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.250"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.217"
$ws.Range("E50").Value = "  +4.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.47"
$ws.Range("E51").Value = "  +0.49%  "
